$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the list date in A1 by one day (17/01/2024 -> 18/01/2024)
$ws.Range("A1").Value = 45309

# Step 2: update the prices for the two items in the list
$ws.Range("D29").Value = 1001
$ws.Range("D30").Value = 1077
